$wb = $excel.ActiveWorkbook

# 1. Duplicate the "2022-Q2" sheet to create a new "2022-Q3" sheet, placed
#    right before it (so the tab order becomes ...,2022-Q3,2022-Q2,...).
$src = $wb.Worksheets.Item("2022-Q2")
$src.Copy($src)
$newws = $wb.Worksheets.Item(2)
$newws.Name = "2022-Q3"

# 2. Update the fund metrics on the new "2022-Q3" sheet with the Q3 figures.
#    Values are kept as text (leading "'") so they stay the same cell type
#    ("0.61" etc.) as the source sheet they were copied from.
$newws.Range("D2").Value = "'0.61"
$newws.Range("E2").Value = "'89.76"
$newws.Range("F2").Value = "'2.97"
$newws.Range("G2").Value = "'0.0181"

# 3. Insert a new summary row for 2022-Q3 at the top of the "总计" sheet,
#    pushing the existing quarters down by one row.
$ws = $wb.Worksheets.Item("总计")
$ws.Rows.Item(2).Insert()

# Copy formatting from the row below (so the new row matches the existing
# look: bold/centered index column, plain data columns) before filling in
# the values.
$ws.Range("A3:D3").Copy()
$ws.Range("A2:D2").PasteSpecial(-4122)

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "2022-Q3"
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 0.02

# Renumber the index column for the rows that shifted down.
$ws.Range("A3").Value = 1
$ws.Range("A4").Value = 2
$ws.Range("A5").Value = 3
$ws.Range("A6").Value = 4

# 4. Restore the originally-active tab ("2020-Q4"), since copying a sheet
#    shifts the active/selected tab as a side effect.
$wb.Worksheets.Item("2020-Q4").Activate()
